$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast error table values for rows 7-11 (Q6-Q9 quarters)
$ws.Range("B7").Value = -0.1268670602450409
$ws.Range("C7").Value = 0.4912998571971119
$ws.Range("D7").Value = 0.4440863547471695
$ws.Range("E7").Value = 0.6663980452756216
$ws.Range("F7").Value = 0.6634903322751415
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = -0.1983642510587261
$ws.Range("C8").Value = 0.4794362171649642
$ws.Range("D8").Value = 0.4333623456850792
$ws.Range("E8").Value = 0.6583026246986102
$ws.Range("F8").Value = 0.6368693867280113
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = -0.2061622739680646
$ws.Range("C9").Value = 0.5255257753659175
$ws.Range("D9").Value = 0.494714894569724
$ws.Range("E9").Value = 0.7033597191833806
$ws.Range("F9").Value = 0.689936695312686
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.08791877842089156
$ws.Range("C10").Value = 0.5688666216736536
$ws.Range("D10").Value = 0.6389144850882613
$ws.Range("E10").Value = 0.7993212652546292
$ws.Range("F10").Value = 0.8269120295493002
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = -0.5284371888516474
$ws.Range("C11").Value = 0.566781001591608
$ws.Range("D11").Value = 0.5781463867675047
$ws.Range("E11").Value = 0.7603593800088908
$ws.Range("F11").Value = 0.6112492578789698
$ws.Range("G11").Value = 5
